$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row indices below are 1-based Word table row numbers
# (0-based row index + 1, as identified from the document structure)

$t.Cell(1, 1).Range.Text = "0M"        # was 99.94
$t.Cell(2, 1).Range.Text = "0M"        # was 0.04
$t.Cell(3, 1).Range.Text = "0M"        # was 55
$t.Cell(4, 1).Range.Text = "218"       # was 78
$t.Cell(5, 1).Range.Text = "0.00001"   # was 0.00005
$t.Cell(7, 1).Range.Text = "0.00014"   # was 0.00012
$t.Cell(9, 1).Range.Text = "0.00022"   # was 0.00014
$t.Cell(10, 1).Range.Text = "0.00025"  # was 0.00015
$t.Cell(11, 1).Range.Text = "0.00031"  # was 0.00018
$t.Cell(12, 1).Range.Text = "0.03556"  # was 0.01228

$t.Cell(44, 1).Range.Text = "99.94"    # was "70\t0.00017\t0.00051\t0.00028\t0.00008\t0.00022\t0.00025\t0.00031\t0.01949\t100.0"
$t.Cell(45, 1).Range.Text = "0.04"     # was "62\t0.00001\t0.00008\t0.00004\t0.00002\t0.00003\t0.00004\t0.00006\t0.00264\t100.0"
$t.Cell(46, 1).Range.Text = "55"       # was "8\t0.00012\t0.00019\t0.00014\t0.00003\t0.00012\t0.00013\t0.00015\t0.00115\t100.0"
